# Generate Report for Handback
#
# Updates the localization-status workbook to reflect a completed handback:
#  - "Ready for handoff" status text becomes "Handed back: in sync with en-US"
#    (Overview!E2:F3 and the Status column on the zh-cn / de-de sheets share
#    this string, so changing the text updates every cell that used it)
#  - The zh-cn and de-de sheets get their "Latest Target File" / "Latest
#    Handback File" columns (I/J) populated, with hyperlinks added on the
#    "Latest Target File" cells
#  - The de-de sheet's "Latest Handback DateTime" (K) gets a real timestamp
#  - The zh-cn sheet's existing (but blank) "Latest Handback DateTime"
#    timestamp placeholder is replaced with a real timestamp as well
#  - A few columns are widened to fit the newly-populated content

$wb = $excel.ActiveWorkbook

$mdFile      = "90ebd6c6-f264-4139-bf5d-d442b8369e6a.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/599b2efc82627c1479b540df59e271317131e0ce/e2e/90ebd6c6-f264-4139-bf5d-d442b8369e6a.md"
$zhXlf       = "90ebd6c6-f264-4139-bf5d-d442b8369e6a.0bd7fec0a99f4912c574be2c26140f9d77143a7d.zh-cn.xlf"
$deXlf       = "90ebd6c6-f264-4139-bf5d-d442b8369e6a.0bd7fec0a99f4912c574be2c26140f9d77143a7d.de-de.xlf"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Every sheet that carries that status string gets refreshed.
# ---------------------------------------------------------------------
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $used = $sheet.UsedRange
    foreach ($cell in $used.Cells) {
        # NB: compare with the literal on the left -- some cells hold a
        # Boolean ("True"/"False"), and "$boolCell.Value2 -eq <string>"
        # would coerce the string to Boolean (any non-empty string is
        # truthy) and false-positive match. Literal-on-the-left makes
        # PowerShell coerce the other way instead.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value2 = "Handed back: in sync with en-US"
        }
    }
}

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I) / Latest Handback File (J)
#    for rows 2 & 3, add hyperlinks on column I, and set the Latest
#    Handback DateTime (K).
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$zh.Range("J2").Value = $zhXlf
$zh.Range("K2").Value = "2016-09-04 09:06:48"

$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$zh.Range("J3").Value = $zhXlf
$zh.Range("K3").Value = "2016-09-04 09:06:48"

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, with the de-de xlf name and its own
#    handback timestamp.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$de.Range("J2").Value = $deXlf
$de.Range("K2").Value = "2016-09-04 09:06:55"

$de.Hyperlinks.Add($de.Range("I3"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFile)
$de.Range("J3").Value = $deXlf
$de.Range("K3").Value = "2016-09-04 09:06:55"

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# 4. Overview sheet: widen the zh-cn / de-de status columns (E, F) to
#    match the longer status text.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527
